$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values regenerated for s_vals data (filtering save games), per row.
# Columns: B=TB, C=d2S, D=K, E=IP, F=Win (unchanged), G=sum (=B+C+D+E)

$data = @{
    2 = @{ B = 3.230985683306322;  C = 1.667794583268128; D = 0.1575252929769615; E = 0.496779210170732;  G = 5.553084769722144  }
    3 = @{ B = 0.003994804209775715; C = 1.667794583268128; D = 26.21740644021617; E = 645.3272768299601; G = 673.2164726576541  }
    4 = @{ B = 3.230985683306322;  C = 1.667794583268128; D = 0.1575252929769615; E = 0.496779210170732;  G = 5.553084769722144  }
    5 = @{ B = 3.230985683306322;  C = 1.667794583268128; D = 0.1575252929769615; E = 0.496779210170732;  G = 5.553084769722144  }
    6 = @{ B = 3.230985683306322;  C = 1.667794583268128; D = 3.900430680208489;  E = 8.660232485948974;  G = 17.45944343273191 }
    7 = @{ B = 0.6753301551942219; C = 1.667794583268128; D = 0.8054896365839992; E = 8.660232485948974;  G = 11.80884686099532 }
    8 = @{ B = 1.459612070389937;  C = 1.667794583268128; D = 3.900430680208489;  E = 0.496779210170732;  G = 7.524616544037286 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals.B   # B
    $ws.Cells.Item($row, 3).Value = $vals.C   # C
    $ws.Cells.Item($row, 4).Value = $vals.D   # D
    $ws.Cells.Item($row, 5).Value = $vals.E   # E
    $ws.Cells.Item($row, 7).Value = $vals.G   # G
}
